$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 3).Value = $excel.Evaluate("2.473e-19")
$ws.Cells.Item(2, 4).Value = 2.429
$ws.Cells.Item(2, 5).Value = 1.715
$ws.Cells.Item(3, 3).Value = $excel.Evaluate("1.81e-18")
$ws.Cells.Item(3, 4).Value = 2.248
$ws.Cells.Item(3, 5).Value = 0.6059
$ws.Cells.Item(4, 3).Value = $excel.Evaluate("5.073e-18")
$ws.Cells.Item(4, 4).Value = 3.451
$ws.Cells.Item(4, 5).Value = 0.8223
$ws.Cells.Item(5, 3).Value = $excel.Evaluate("6.555e-19")
$ws.Cells.Item(5, 4).Value = 2.633
$ws.Cells.Item(5, 5).Value = 0.1534
$ws.Cells.Item(6, 3).Value = $excel.Evaluate("1.301e-18")
$ws.Cells.Item(6, 4).Value = 2.951
$ws.Cells.Item(6, 5).Value = 0.1706
$ws.Cells.Item(7, 3).Value = 0.351
$ws.Cells.Item(7, 4).Value = 1.446
$ws.Cells.Item(7, 5).Value = 0.4955
$ws.Cells.Item(7, 6).Value = 25
$ws.Cells.Item(8, 3).Value = $excel.Evaluate("3.792e-17")
$ws.Cells.Item(8, 4).Value = 3.122
$ws.Cells.Item(8, 5).Value = 2.671
$ws.Cells.Item(9, 3).Value = $excel.Evaluate("9.334e-18")
$ws.Cells.Item(9, 4).Value = 5.336
$ws.Cells.Item(9, 5).Value = 2.182
$ws.Cells.Item(10, 3).Value = $excel.Evaluate("4.402e-18")
$ws.Cells.Item(10, 4).Value = 3.549
$ws.Cells.Item(10, 5).Value = 0.5923
$ws.Cells.Item(11, 3).Value = 0.09017
$ws.Cells.Item(11, 4).Value = 2.41
$ws.Cells.Item(11, 5).Value = 0.1071
$ws.Cells.Item(11, 6).Value = 10
$ws.Cells.Item(12, 3).Value = $excel.Evaluate("6.26e-19")
$ws.Cells.Item(12, 4).Value = 3.238
$ws.Cells.Item(12, 5).Value = 0.7953
$ws.Cells.Item(13, 3).Value = $excel.Evaluate("2.761e-19")
$ws.Cells.Item(13, 4).Value = 3.334
$ws.Cells.Item(13, 5).Value = 0.9553
$ws.Cells.Item(14, 3).Value = $excel.Evaluate("3.738e-18")
$ws.Cells.Item(14, 4).Value = 3.952
$ws.Cells.Item(14, 5).Value = 1.361
$ws.Cells.Item(14, 6).Value = 5
$ws.Cells.Item(15, 3).Value = $excel.Evaluate("5.053e-18")
$ws.Cells.Item(15, 4).Value = 4.032
$ws.Cells.Item(15, 5).Value = 1.121
$ws.Cells.Item(15, 6).Value = 5
$ws.Cells.Item(16, 3).Value = $excel.Evaluate("4.052e-18")
$ws.Cells.Item(16, 4).Value = 4.107
$ws.Cells.Item(16, 5).Value = 1.278
$ws.Cells.Item(17, 3).Value = $excel.Evaluate("4.711e-18")
$ws.Cells.Item(17, 4).Value = 4.491
$ws.Cells.Item(17, 5).Value = 1.675
$ws.Cells.Item(18, 3).Value = $excel.Evaluate("2.61e-19")
$ws.Cells.Item(18, 4).Value = 4.287
$ws.Cells.Item(18, 5).Value = 0.9827
$ws.Cells.Item(19, 3).Value = 0.003197
$ws.Cells.Item(19, 4).Value = 3.884
$ws.Cells.Item(19, 5).Value = 0.2202
$ws.Cells.Item(19, 6).Value = 15
$ws.Cells.Item(20, 3).Value = $excel.Evaluate("1.423e-17")
$ws.Cells.Item(20, 4).Value = 6.008
$ws.Cells.Item(20, 5).Value = 0.2463
$ws.Cells.Item(21, 3).Value = $excel.Evaluate("1.188e-18")
$ws.Cells.Item(21, 4).Value = 5.064
$ws.Cells.Item(21, 5).Value = 0.3729
$ws.Cells.Item(22, 3).Value = $excel.Evaluate("3.725e-19")
$ws.Cells.Item(22, 4).Value = 5.486
$ws.Cells.Item(22, 5).Value = 0.3875
$ws.Cells.Item(23, 3).Value = $excel.Evaluate("2.309e-19")
$ws.Cells.Item(23, 4).Value = 5.572
$ws.Cells.Item(23, 5).Value = 0.4257
$ws.Cells.Item(24, 3).Value = $excel.Evaluate("5.297e-19")
$ws.Cells.Item(24, 4).Value = 6.303
$ws.Cells.Item(24, 5).Value = 0.6127
$ws.Cells.Item(25, 3).Value = $excel.Evaluate("5.022e-20")
$ws.Cells.Item(25, 4).Value = 5.954
$ws.Cells.Item(25, 5).Value = 0.841
$ws.Cells.Item(25, 6).Value = 5
$ws.Cells.Item(26, 3).Value = 0.001364
$ws.Cells.Item(26, 4).Value = 4.536
$ws.Cells.Item(26, 5).Value = 0.3336
$ws.Cells.Item(26, 6).Value = 10
$ws.Cells.Item(27, 3).Value = $excel.Evaluate("1.771e-18")
$ws.Cells.Item(27, 4).Value = 5.316
$ws.Cells.Item(27, 5).Value = 0.5809
$ws.Cells.Item(28, 3).Value = $excel.Evaluate("4.031e-19")
$ws.Cells.Item(28, 4).Value = 6.067
$ws.Cells.Item(28, 5).Value = 0.7125
$ws.Cells.Item(29, 3).Value = $excel.Evaluate("4.213e-21")
$ws.Cells.Item(29, 4).Value = 6.252
$ws.Cells.Item(29, 5).Value = 0.4791
$ws.Cells.Item(30, 3).Value = $excel.Evaluate("2.488e-21")
$ws.Cells.Item(30, 4).Value = 6.048
$ws.Cells.Item(30, 5).Value = 0.3386
$ws.Cells.Item(31, 3).Value = 0.0002918
$ws.Cells.Item(31, 4).Value = 5.026
$ws.Cells.Item(31, 5).Value = 0.253
$ws.Cells.Item(31, 6).Value = 10
$ws.Cells.Item(32, 3).Value = $excel.Evaluate("6.51e-19")
$ws.Cells.Item(32, 4).Value = 5.849
$ws.Cells.Item(32, 5).Value = 0.922
$ws.Cells.Item(33, 3).Value = $excel.Evaluate("5.907000000000001e-20")
$ws.Cells.Item(33, 4).Value = 5.812
$ws.Cells.Item(33, 5).Value = 0.7369
$ws.Cells.Item(34, 3).Value = $excel.Evaluate("4.038e-18")
$ws.Cells.Item(34, 4).Value = 6.829
$ws.Cells.Item(34, 5).Value = 0.9738
$ws.Cells.Item(34, 6).Value = 5
$ws.Cells.Item(35, 3).Value = $excel.Evaluate("1.077e-18")
$ws.Cells.Item(35, 4).Value = 6.376
$ws.Cells.Item(35, 5).Value = 1.074
$ws.Cells.Item(36, 3).Value = $excel.Evaluate("4.016e-19")
$ws.Cells.Item(36, 4).Value = 6.471
$ws.Cells.Item(36, 5).Value = 1.065
$ws.Cells.Item(37, 3).Value = $excel.Evaluate("1.252e-18")
$ws.Cells.Item(37, 4).Value = 6.166
$ws.Cells.Item(37, 5).Value = 0.2847
$ws.Cells.Item(38, 3).Value = $excel.Evaluate("2.446e-19")
$ws.Cells.Item(38, 4).Value = 6.72
$ws.Cells.Item(38, 5).Value = 0.4637
$ws.Cells.Item(39, 3).Value = $excel.Evaluate("4.01e-21")
$ws.Cells.Item(39, 4).Value = 6.863
$ws.Cells.Item(39, 5).Value = 0.7395
$ws.Cells.Item(40, 3).Value = $excel.Evaluate("1.722e-19")
$ws.Cells.Item(40, 4).Value = 6.569
$ws.Cells.Item(40, 5).Value = 0.8516
$ws.Cells.Item(41, 3).Value = $excel.Evaluate("1.741e-18")
$ws.Cells.Item(41, 4).Value = 6.26
$ws.Cells.Item(41, 5).Value = 1.042
$ws.Cells.Item(41, 6).Value = 5
$ws.Cells.Item(42, 3).Value = $excel.Evaluate("1.384e-19")
$ws.Cells.Item(42, 4).Value = 6.623
$ws.Cells.Item(42, 5).Value = 1.246
$ws.Cells.Item(43, 3).Value = 0.06199
$ws.Cells.Item(43, 4).Value = 4.316
$ws.Cells.Item(43, 5).Value = 0.2819
$ws.Cells.Item(43, 6).Value = 20
$ws.Cells.Item(44, 3).Value = $excel.Evaluate("2.274e-17")
$ws.Cells.Item(44, 4).Value = 6.905
$ws.Cells.Item(44, 5).Value = 1.989
$ws.Cells.Item(45, 3).Value = $excel.Evaluate("6.644e-19")
$ws.Cells.Item(45, 4).Value = 6.15
$ws.Cells.Item(45, 5).Value = 2.506
$ws.Cells.Item(46, 3).Value = $excel.Evaluate("1.535e-19")
$ws.Cells.Item(46, 4).Value = 6.803
$ws.Cells.Item(46, 5).Value = 2.1
$ws.Cells.Item(47, 3).Value = $excel.Evaluate("8.204e-19")
$ws.Cells.Item(47, 4).Value = 7.293
$ws.Cells.Item(47, 5).Value = 0.8519
$ws.Cells.Item(48, 3).Value = 0.406
$ws.Cells.Item(48, 4).Value = 6.291
$ws.Cells.Item(48, 5).Value = 0.144
$ws.Cells.Item(48, 6).Value = 25
$ws.Cells.Item(49, 3).Value = 0.003674
$ws.Cells.Item(49, 4).Value = 5.464
$ws.Cells.Item(49, 5).Value = 0.7191
$ws.Cells.Item(49, 6).Value = 10
$ws.Cells.Item(50, 3).Value = $excel.Evaluate("2.14e-17")
$ws.Cells.Item(50, 4).Value = 5.941
$ws.Cells.Item(50, 5).Value = 2.329
$ws.Cells.Item(50, 6).Value = 5
$ws.Cells.Item(51, 3).Value = $excel.Evaluate("1.261e-17")
$ws.Cells.Item(51, 4).Value = 6.84
$ws.Cells.Item(51, 5).Value = 3.537
$ws.Cells.Item(52, 3).Value = $excel.Evaluate("1.349e-18")
$ws.Cells.Item(52, 4).Value = 7.2
$ws.Cells.Item(52, 5).Value = 3.006
$ws.Cells.Item(53, 3).Value = $excel.Evaluate("1.237e-18")
$ws.Cells.Item(53, 4).Value = 6.375
$ws.Cells.Item(53, 5).Value = 2.212
$ws.Cells.Item(54, 3).Value = $excel.Evaluate("1.935e-18")
$ws.Cells.Item(54, 4).Value = 6.707
$ws.Cells.Item(54, 5).Value = 0.4235
$ws.Cells.Item(55, 3).Value = $excel.Evaluate("5.889e-19")
$ws.Cells.Item(55, 4).Value = 5.466
$ws.Cells.Item(55, 5).Value = 1.096
$ws.Cells.Item(56, 3).Value = $excel.Evaluate("2.337e-18")
$ws.Cells.Item(56, 4).Value = 5.359
$ws.Cells.Item(56, 5).Value = 1.331
$ws.Cells.Item(57, 3).Value = $excel.Evaluate("2.15e-18")
$ws.Cells.Item(57, 4).Value = 6.439
$ws.Cells.Item(57, 5).Value = 1.968
$ws.Cells.Item(58, 3).Value = $excel.Evaluate("3.592e-18")
$ws.Cells.Item(58, 4).Value = 4.914
$ws.Cells.Item(58, 5).Value = 1.268
$ws.Cells.Item(59, 3).Value = $excel.Evaluate("1.538e-18")
$ws.Cells.Item(59, 4).Value = 5.414
$ws.Cells.Item(59, 5).Value = 0.3571
$ws.Cells.Item(60, 3).Value = $excel.Evaluate("8.575e-19")
$ws.Cells.Item(60, 4).Value = 6.023
$ws.Cells.Item(60, 5).Value = 0.2395
$ws.Cells.Item(61, 4).Value = 6.576
$ws.Cells.Item(61, 5).Value = 0.1818
$ws.Cells.Item(61, 6).Value = 25
$ws.Cells.Item(62, 3).Value = 0.147
$ws.Cells.Item(62, 4).Value = 6.771
$ws.Cells.Item(62, 5).Value = 0.1987
$ws.Cells.Item(62, 6).Value = 25
$ws.Cells.Item(63, 3).Value = $excel.Evaluate("1.016e-17")
$ws.Cells.Item(63, 4).Value = 6.383
$ws.Cells.Item(63, 5).Value = 0.9268999999999999
$ws.Cells.Item(64, 3).Value = $excel.Evaluate("8.577e-19")
$ws.Cells.Item(64, 4).Value = 5.523
$ws.Cells.Item(64, 5).Value = 1.142
$ws.Cells.Item(65, 3).Value = $excel.Evaluate("6.728e-18")
$ws.Cells.Item(65, 4).Value = 6.621
$ws.Cells.Item(65, 5).Value = 1.437
$ws.Cells.Item(66, 3).Value = $excel.Evaluate("1.294e-17")
$ws.Cells.Item(66, 4).Value = 7.142
$ws.Cells.Item(66, 5).Value = 2.165
$ws.Cells.Item(67, 3).Value = $excel.Evaluate("7.243e-19")
$ws.Cells.Item(67, 4).Value = 6.465
$ws.Cells.Item(67, 5).Value = 1.611
$ws.Cells.Item(68, 3).Value = $excel.Evaluate("2.539e-19")
$ws.Cells.Item(68, 4).Value = 7.099
$ws.Cells.Item(68, 5).Value = 0.5243
$ws.Cells.Item(69, 3).Value = 0.006925
$ws.Cells.Item(69, 4).Value = 4.937
$ws.Cells.Item(69, 5).Value = 0.2123
$ws.Cells.Item(69, 6).Value = 10
$ws.Cells.Item(70, 3).Value = 0.0007808
$ws.Cells.Item(70, 4).Value = 4.815
$ws.Cells.Item(70, 5).Value = 0.3495
$ws.Cells.Item(70, 6).Value = 10
$ws.Cells.Item(71, 3).Value = $excel.Evaluate("1.359e-19")
$ws.Cells.Item(71, 4).Value = 6.906
$ws.Cells.Item(71, 5).Value = 1.134
$ws.Cells.Item(72, 3).Value = $excel.Evaluate("7.207e-19")
$ws.Cells.Item(72, 4).Value = 8.228
$ws.Cells.Item(72, 5).Value = 0.9114
$ws.Cells.Item(72, 6).Value = 5
